# Mise au propre d'une grande partie du code et clarification a quelques endroit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw input values (H15, H16, H17)
$ws.Range("H15").Value = 90
$ws.Range("H16").Value = 160
$ws.Range("H17").Value = 0

# Update the view: scroll position + active selection
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("H16").Select()
